$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the old totals row (old row 29 -> becomes row 30).
$ws.Rows.Item(29).Insert()

# The freshly inserted row 29 should look like the other data rows (5-28):
# copy formatting (styles/borders) down from row 28.
$ws.Range("A28:I28").Copy()
$ws.Range("A29:I29").PasteSpecial(-4122)

# Fix up the totals row (now row 30): extend the SUM ranges to include the
# new row 29, and add the missing "Da ho tro" (G) total.
$ws.Range("E30").Formula = "=SUM(E`$5:E29)"
$ws.Range("G30").Formula = "=SUM(G`$5:G29)"

# Restore the current selection to B9.
[void]$ws.Range("B9").Select()
